$wb = $excel.ActiveWorkbook

# --- Sheet "Summary": remove leftover formatting rows 8:13 (no real data), shrinking used range to A1:F5 ---
$wsSummary = $wb.Worksheets.Item("Summary")
$wsSummary.Rows("8:13").Delete()
$wsSummary.Range("C4").Select()

# --- Sheet "Repayment Schedule": move the "Over Due" column (header + data rows 3:14) from O to P ---
$wsSchedule = $wb.Worksheets.Item("Repayment Schedule")
$wsSchedule.Range("O1").Cut($wsSchedule.Range("P1"))
$wsSchedule.Range("O3:O14").Cut($wsSchedule.Range("P3:P14"))
$wsSchedule.Range("G11").Select()

# --- Sheet "Transactions": remove leftover / stray formatting cells and rows outside the real data area ---
$wsTransactions = $wb.Worksheets.Item("Transactions")
$wsTransactions.Range("K2:M6").Clear()
$wsTransactions.Range("A5:D6").Clear()
$wsTransactions.Range("E21").Clear()
$wsTransactions.Activate()
$wsTransactions.Range("A1").Select()
